# Apply updated crypto market data values while preserving the
# original plain-text (inline string) representation of each cell.
# Columns D/E hold numeric-looking text (e.g. "248.30", "0.664") that
# Excel would otherwise silently reinterpret as real numbers, so for
# column D we briefly force a Text number format, assign the literal
# string, then restore the cell to the workbook's default "Normal"
# style (no explicit style index), matching the original formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.097.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.045.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.664"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.383"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0785"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.331.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.832"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.046.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "17.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +22.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.107.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "75.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0895"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("E26").Value = "  +7.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.47%  "
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0623"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.15%  "
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("E36").Value = "  -3.07%  "
$ws.Range("E37").Value = "  -1.42%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("B39").Value = "Cronos"
$ws.Range("C39").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.106"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.67%  "
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.10%  "
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.46%  "
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.283.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.233.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -21.22%  "
